$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7518376611942944
$ws.Range("B3").Value = 0.7446116706216693
$ws.Range("B4").Value = 0.8909549322393191
$ws.Range("B5").Value = 0.7482067196427615
$ws.Range("B6").Value = 0.7168857455253601
$ws.Range("B7").Value = 0.6202630996704102
